# Apply updated crypto price / volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.168.66"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "1.905.93"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'326.32"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.4609"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.3891"
$ws.Range("E8").Value = "  -1.25%  "

$ws.Range("D9").Value = "'0.07881"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").Value = "'0.9910"
$ws.Range("E10").Value = "  -1.25%  "

$ws.Range("D11").Value = "'21.94"
$ws.Range("E11").Value = "  -1.80%  "

$ws.Range("D12").Value = "1.908.96"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "'5.774"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("D14").Value = "'7.049"

$ws.Range("D15").Value = "'0.07030"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").Value = "'88.09"
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "'0.000009951"
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "29.171.84"

$ws.Range("D22").Value = "'5.322"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("D23").Value = "'11.15"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").Value = "'2.104"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").Value = "'155.90"
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").Value = "'19.42"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("D27").Value = "'5.931"
$ws.Range("E27").Value = "  -3.18%  "

$ws.Range("D28").Value = "'118.72"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").Value = "'1.885"
$ws.Range("E29").Value = "  -5.59%  "

$ws.Range("D30").Value = "'0.09355"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").Value = "'0.8960"
$ws.Range("E31").Value = "  -3.31%  "

$ws.Range("D32").Value = "'5.244"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("D33").Value = "'1.320"
$ws.Range("E33").Value = "  -2.55%  "

$ws.Range("D34").Value = "'3.161"
$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("E35").Value = "  -0.65%  "

$ws.Range("D36").Value = "'1.172"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("D37").Value = "'0.02089"
$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").Value = "'7.688"
$ws.Range("E39").Value = "  -3.39%  "

$ws.Range("D40").Value = "'0.5707"
$ws.Range("E40").Value = "  -0.85%  "

$ws.Range("D41").Value = "'0.1797"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").Value = "'9.729"
$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("D43").Value = "'11.90"
$ws.Range("E43").Value = "  -1.14%  "

$ws.Range("D44").Value = "'0.5360"
$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("D45").Value = "'2.181"
$ws.Range("E45").Value = "  -3.83%  "

$ws.Range("D46").Value = "'0.07021"
$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("D47").Value = "'1.849"
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("D48").Value = "'2.549"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").Value = "'113.11"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("D50").Value = "'0.2945"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'71.42"
$ws.Range("E51").Value = "  -0.37%  "
